# Logged Week 15 and simulated Week 16
# Update Row 3 ("R" / Road) target-depth stats on both the OFF and DEF sheets.

$wb = $excel.ActiveWorkbook

# OFF sheet: Short Att, Short Comp, Deep Att, Short Int for the "R" row (row 3)
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 191
$wsOff.Range("C3").Value = 135
$wsOff.Range("D3").Value = 34
$wsOff.Range("F3").Value = 4

# DEF sheet: Short Att, Short Comp, Deep Att, Deep Int for the "R" row (row 3)
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 221
$wsDef.Range("C3").Value = 151
$wsDef.Range("D3").Value = 53
$wsDef.Range("G3").Value = 3
